$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.017009496688843
$ws.Range("B1").Value = 2.060846090316772
$ws.Range("C1").Value = 3.874757766723633
$ws.Range("D1").Value = 1.273569941520691
$ws.Range("E1").Value = 0.6976171135902405
